$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text format on cells whose new values look like pure numbers,
# so Excel stores them as text (matching the source inlineStr cells)
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D50').NumberFormat = '@'

# Apply cell value updates
$ws.Range('D2').Value() = '54.606.43'
$ws.Range('E2').Value() = '  -5.33%  '
$ws.Range('D3').Value() = '2.458.46'
$ws.Range('E3').Value() = '  -8.26%  '
$ws.Range('D4').Value() = '1.00'
$ws.Range('E4').Value() = '  -0.06%  '
$ws.Range('D5').Value() = '470.24'
$ws.Range('E5').Value() = '  -4.65%  '
$ws.Range('D6').Value() = '134.26'
$ws.Range('E6').Value() = '  +0.18%  '
$ws.Range('D7').Value() = '0.996'
$ws.Range('E7').Value() = '  -0.86%  '
$ws.Range('D8').Value() = '0.493'
$ws.Range('E8').Value() = '  -4.95%  '
$ws.Range('D9').Value() = '2.474.77'
$ws.Range('E9').Value() = '  -6.65%  '
$ws.Range('D10').Value() = '0.0967'
$ws.Range('E10').Value() = '  -3.96%  '
$ws.Range('D11').Value() = '5.35'
$ws.Range('E11').Value() = '  -7.57%  '
$ws.Range('D12').Value() = '0.320'
$ws.Range('E12').Value() = '  -5.17%  '
$ws.Range('E13').Value() = '  -3.18%  '
$ws.Range('D14').Value() = '2.893.52'
$ws.Range('E14').Value() = '  -8.91%  '
$ws.Range('D15').Value() = '54.478.84'
$ws.Range('E15').Value() = '  -5.83%  '
$ws.Range('B16').Value() = 'Avalanche'
$ws.Range('C16').Value() = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$ws.Range('D16').Value() = '20.03'
$ws.Range('E16').Value() = '  -4.14%  '
$ws.Range('B17').Value() = 'ShibaInu'
$ws.Range('C17').Value() = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range('D17').Value() = '0.0000133'
$ws.Range('E17').Value() = '  +1.14%  '
$ws.Range('D18').Value() = '2.469.57'
$ws.Range('E18').Value() = '  -8.63%  '
$ws.Range('D19').Value() = '4.24'
$ws.Range('E19').Value() = '  -7.65%  '
$ws.Range('D20').Value() = '311.53'
$ws.Range('E20').Value() = '  -7.76%  '
$ws.Range('D21').Value() = '9.49'
$ws.Range('E21').Value() = '  -10.56%  '
$ws.Range('E22').Value() = '  +0.59%  '
$ws.Range('D23').Value() = '5.69'
$ws.Range('E23').Value() = '  +1.70%  '
$ws.Range('D24').Value() = '5.40'
$ws.Range('E24').Value() = '  -10.96%  '
$ws.Range('D25').Value() = '57.15'
$ws.Range('E25').Value() = '  -6.59%  '
$ws.Range('E26').Value() = '  +0.80%  '
$ws.Range('D27').Value() = '0.389'
$ws.Range('E27').Value() = '  -6.12%  '
$ws.Range('D28').Value() = '2.549.09'
$ws.Range('E28').Value() = '  -10.43%  '
$ws.Range('D29').Value() = '0.155'
$ws.Range('E29').Value() = '  -7.46%  '
$ws.Range('D30').Value() = '7.30'
$ws.Range('E30').Value() = '  +1.14%  '
$ws.Range('D31').Value() = '0.999'
$ws.Range('E31').Value() = '  -0.28%  '
$ws.Range('D32').Value() = '0.0₃0734'
$ws.Range('E32').Value() = '  -7.08%  '
$ws.Range('D33').Value() = '150.26'
$ws.Range('E33').Value() = '  +2.60%  '
$ws.Range('D34').Value() = '17.95'
$ws.Range('E34').Value() = '  -3.20%  '
$ws.Range('D35').Value() = '1.45'
$ws.Range('E35').Value() = '  -7.44%  '
$ws.Range('D36').Value() = '5.07'
$ws.Range('E36').Value() = '  -1.88%  '
$ws.Range('D37').Value() = '3.59'
$ws.Range('E37').Value() = '  -11.72%  '
$ws.Range('E38').Value() = '  -2.50%  '
$ws.Range('D39').Value() = '0.813'
$ws.Range('E39').Value() = '  -9.14%  '
$ws.Range('D40').Value() = '33.87'
$ws.Range('E40').Value() = '  -5.83%  '
$ws.Range('D41').Value() = '0.992'
$ws.Range('E41').Value() = '  -0.38%  '
$ws.Range('D42').Value() = '0.608'
$ws.Range('E42').Value() = '  +4.25%  '
$ws.Range('D43').Value() = '0.0534'
$ws.Range('E43').Value() = '  -1.48%  '
$ws.Range('D44').Value() = '3.31'
$ws.Range('E44').Value() = '  -3.42%  '
$ws.Range('B45').Value() = 'Stacks'
$ws.Range('C45').Value() = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range('D45').Value() = '1.26'
$ws.Range('E45').Value() = '  -5.16%  '
$ws.Range('B46').Value() = 'WhiteBITCoin'
$ws.Range('C46').Value() = 'https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt'
$ws.Range('D46').Value() = '10.20'
$ws.Range('E46').Value() = '  -1.45%  '
$ws.Range('D47').Value() = '1.960.96'
$ws.Range('E47').Value() = '  -7.87%  '
$ws.Range('D48').Value() = '0.0221'
$ws.Range('E48').Value() = '  -0.40%  '
$ws.Range('D49').Value() = '0.0879'
$ws.Range('E49').Value() = '  +0.84%  '
$ws.Range('D50').Value() = '4.31'
$ws.Range('E50').Value() = '  -4.57%  '
$ws.Range('E51').Value() = '  -9.10%  '
